# Generate Report for Handoff
#
# The "b47cd761-22d4-43b9-9d16-7370abf34e62" row's handoff run completed at a
# later timestamp, and that newer "Latest Handoff Datetime" is the value that
# should now be reflected everywhere it had previously been stamped with a
# stale value (rows 12/13 of each language sheet held an older, now-obsolete
# timestamp that needs to collapse onto the current one).

$wb = $excel.ActiveWorkbook

$updates = @{
    "zh-cn" = "2016-03-11 00:34:29"
    "de-de" = "2016-03-11 00:34:37"
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $newTimestamp = $updates[$sheetName]

    # D7  - b47cd761-...: the handoff run that produced the new timestamp
    $ws.Range("D7").Value = $newTimestamp
    # D12 - 40b1cf56-...: previously stamped with the now-stale timestamp
    $ws.Range("D12").Value = $newTimestamp
    # D13 - 8a1371c3-...: previously stamped with the now-stale timestamp
    $ws.Range("D13").Value = $newTimestamp
}
